$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.447.20"
$ws.Range("E2").Value = "  -0.28%  "

$ws.Range("D3").Value = "1.727.34"
$ws.Range("E3").Value = "  -0.57%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9974"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.80%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9978"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.18%  "

$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2604"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.55%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06200"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.61%  "

$ws.Range("D10").Value = "1.728.90"
$ws.Range("E10").Value = "  -0.40%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06997"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.50%  "

$ws.Range("E12").Value = "  -0.69%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.529"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.55%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6010"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.58%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.14"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.44%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9975"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.21%  "

$ws.Range("D17").Value = "26.445.33"
$ws.Range("E17").Value = "  -0.26%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9973"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.25%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007165"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.00%  "

$ws.Range("E20").Value = "  -1.73%  "

$ws.Range("D21").Value = "1.945.98"
$ws.Range("E21").Value = "  -0.51%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.470"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.85%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.506"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.106"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.58%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "137.52"
$ws.Range("D25").Style = "Normal"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.39%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.413"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.04%  "

$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.750"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.56%  "

$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "106.64"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.32%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.910"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.72%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08016"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.74%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.642"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.43%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04493"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.70%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9968"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.26%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.600"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.35%  "

$ws.Range("E36").Value = "  -0.95%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6233"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.00%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9393"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.93%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.992"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.09%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.386"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.48%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9973"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.51%  "

$ws.Range("E42").Value = "  -1.82%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.86"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.41%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.410"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.42%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3851"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.24%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.913"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.48%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1159"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.24%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05370"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.38%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.51"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.09%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.747"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.35%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.39"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.66%  "
